$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-11-24"

# Update the column header label (shared string) used by cell I1
$ws.Range("I1").Value = "2022 (through 11-24)"

# Update November value for 2022 (through 11-24) column
$ws.Range("I12").Value = 88

# Update Total value for 2022 (through 11-24) column
$ws.Range("I14").Value = 1486
